# Servo regulator PTH resistors: switch from 1/6W PTH resistor stock (no
# vendor data) to vertical-mount 1/4W (0207/5V) Digi-Key parts, since the
# needed resistance values aren't available in 1/6W PTH and a normal 1/4W
# won't fit lengthwise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R37 (1.96K) -> row 55
$ws.Range("C55").Value = "R-US_0207/5V"
$ws.Range("D55").Value = "0207/5V"
$ws.Range("E55").Value = "DK"
$ws.Range("F55").Value = "RNF14FTD1K96CT-ND"
$ws.Range("G55").Value = "RNF14FTD1K96"
$ws.Range("I55").Value = 0.15

# R38 (221K) -> row 56
$ws.Range("C56").Value = "R-US_0207/5V"
$ws.Range("D56").Value = "0207/5V"
$ws.Range("E56").Value = "DK"
$ws.Range("F56").Value = "RNF14FTD221KCT-ND"
$ws.Range("G56").Value = "RNF14FTD221K"
$ws.Range("I56").Value = 0.15

# R145 (41.2K) -> row 72
$ws.Range("C72").Value = "R-US_0207/5V"
$ws.Range("D72").Value = "0207/5V"
$ws.Range("F72").Value = "RNF14FTD41K2CT-ND"
$ws.Range("G72").Value = "RNF14FTD41K2"
$ws.Range("I72").Value = 0.15

# R146 (66.5K) -> row 73
$ws.Range("C73").Value = "R-US_0207/5V"
$ws.Range("D73").Value = "0207/5V"
$ws.Range("F73").Value = "RNF14FTD66K5CT-ND"
$ws.Range("G73").Value = "RNF14FTD66K5"
$ws.Range("I73").Value = 0.15

# Move the saved view/selection to match the author's final cursor position.
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("I55").Select()
